# Insert a new weekly price block (4 rows: Especial/Primera/Segunda/Tercera)
# for Piña - Vega Central Mapocho de Santiago at the top of the existing
# data table (row 676), shifting every subsequent row down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 676:782 down to 680:786 (inserts 4 blank rows at 676; all
# data/formatting below moves down automatically).
$ws.Range("A676:T679").EntireRow.Insert()

# Populate the newly-inserted rows with the new week's data
# (fecha = 2021-11-05, serial 44505).

# Row 676 - Especial
$ws.Range("A676").Value = 9
$ws.Range("B676").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C676").Value = 'Metropolitana'
$ws.Range("D676").Value = 44505
$ws.Range("E676").Value = 13
$ws.Range("F676").Value = 'Fruta'
$ws.Range("G676").Value = 100108
$ws.Range("H676").Value = 'Tropicales y subtropicales'
$ws.Range("I676").Value = 100108005
$ws.Range("J676").Value = 'Piña'
$ws.Range("K676").Value = 'Caramelo'
$ws.Range("L676").Value = 'Especial'
$ws.Range("M676").Value = 25
$ws.Range("N676").Value = 17000
$ws.Range("O676").Value = 18000
$ws.Range("P676").Value = 17600
$ws.Range("Q676").Value = '$/caja 10 unidades'
$ws.Range("R676").Value = 'Ecuador'
$ws.Range("S676").Value = 1760
$ws.Range("T676").Value = 10

# Row 677 - Primera
$ws.Range("A677").Value = 9
$ws.Range("B677").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C677").Value = 'Metropolitana'
$ws.Range("D677").Value = 44505
$ws.Range("E677").Value = 13
$ws.Range("F677").Value = 'Fruta'
$ws.Range("G677").Value = 100108
$ws.Range("H677").Value = 'Tropicales y subtropicales'
$ws.Range("I677").Value = 100108005
$ws.Range("J677").Value = 'Piña'
$ws.Range("K677").Value = 'Caramelo'
$ws.Range("L677").Value = 'Primera'
$ws.Range("M677").Value = 35
$ws.Range("N677").Value = 17000
$ws.Range("O677").Value = 18000
$ws.Range("P677").Value = 17429
$ws.Range("Q677").Value = '$/caja 12 unidades'
$ws.Range("R677").Value = 'Ecuador'
$ws.Range("S677").Value = 1452
$ws.Range("T677").Value = 12

# Row 678 - Segunda
$ws.Range("A678").Value = 9
$ws.Range("B678").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C678").Value = 'Metropolitana'
$ws.Range("D678").Value = 44505
$ws.Range("E678").Value = 13
$ws.Range("F678").Value = 'Fruta'
$ws.Range("G678").Value = 100108
$ws.Range("H678").Value = 'Tropicales y subtropicales'
$ws.Range("I678").Value = 100108005
$ws.Range("J678").Value = 'Piña'
$ws.Range("K678").Value = 'Caramelo'
$ws.Range("L678").Value = 'Segunda'
$ws.Range("M678").Value = 30
$ws.Range("N678").Value = 17000
$ws.Range("O678").Value = 18000
$ws.Range("P678").Value = 17500
$ws.Range("Q678").Value = '$/caja 14 unidades'
$ws.Range("R678").Value = 'Ecuador'
$ws.Range("S678").Value = 1250
$ws.Range("T678").Value = 14

# Row 679 - Tercera
$ws.Range("A679").Value = 9
$ws.Range("B679").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C679").Value = 'Metropolitana'
$ws.Range("D679").Value = 44505
$ws.Range("E679").Value = 13
$ws.Range("F679").Value = 'Fruta'
$ws.Range("G679").Value = 100108
$ws.Range("H679").Value = 'Tropicales y subtropicales'
$ws.Range("I679").Value = 100108005
$ws.Range("J679").Value = 'Piña'
$ws.Range("K679").Value = 'Caramelo'
$ws.Range("L679").Value = 'Tercera'
$ws.Range("M679").Value = 30
$ws.Range("N679").Value = 17000
$ws.Range("O679").Value = 18000
$ws.Range("P679").Value = 17667
$ws.Range("Q679").Value = '$/caja 16 unidades'
$ws.Range("R679").Value = 'Ecuador'
$ws.Range("S679").Value = 1104
$ws.Range("T679").Value = 16
